# Refresh the cryptos list: update the "Price" (column D) and "Volume(1h)"
# (column E) figures for the rows whose values changed in this data pull.
#
# Column D sometimes holds values that *look* numeric (e.g. "324.21"). If we
# just assign .Value, Excel reinterprets them as real numbers, which loses
# the exact original text (e.g. "41.51" becomes a binary float like
# 41.509999999999998) and changes the cell type away from text. Since the
# source workbook stores these as plain text, we force text interpretation
# (NumberFormat "@") for any replacement value that parses as a plain
# decimal number, then restore the cell to the default "Normal" style so no
# stray number formatting is left behind. Column E values (e.g.
# "  -0.22%  ") never parse as numbers, so they're safe to assign directly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($Cell, $Value) {
    $range = $ws.Range($Cell)

    # Plain decimal number (optional sign, digits, optional single decimal
    # point) -- the kind of text Excel would auto-convert to a numeric cell.
    $looksNumeric = $Value -match '^[+-]?\d+(\.\d+)?$'

    if ($looksNumeric) {
        $range.NumberFormat = "@"
        $range.Value = $Value
        $range.Style = "Normal"
    } else {
        $range.Value = $Value
    }
}

Set-TextValue "D2"  "27.574.71"
Set-TextValue "E2"  "  -0.22%  "
Set-TextValue "D3"  "1.753.04"
Set-TextValue "E4"  "  +0.04%  "
Set-TextValue "D5"  "324.21"
Set-TextValue "E5"  "  -0.10%  "
Set-TextValue "E6"  "  +0.08%  "
Set-TextValue "D7"  "0.4484"
Set-TextValue "E7"  "  +4.20%  "
Set-TextValue "D8"  "0.3552"
Set-TextValue "E8"  "  -1.48%  "
Set-TextValue "D9"  "0.07454"
Set-TextValue "E9"  "  -1.54%  "
Set-TextValue "D10" "41.51"
Set-TextValue "E10" "  -1.42%  "
Set-TextValue "D11" "1.081"
Set-TextValue "E11" "  -2.61%  "
Set-TextValue "E12" "  +0.06%  "
Set-TextValue "D13" "20.73"
Set-TextValue "E13" "  -0.33%  "
Set-TextValue "D14" "5.983"
Set-TextValue "E14" "  -1.46%  "
Set-TextValue "D15" "7.141"
Set-TextValue "E15" "  -1.25%  "
Set-TextValue "D16" "1.759.50"
Set-TextValue "E16" "  -0.07%  "
Set-TextValue "D17" "93.50"
Set-TextValue "E17" "  +1.39%  "
Set-TextValue "E18" "  -1.08%  "
Set-TextValue "D19" "0.06458"
Set-TextValue "E19" "  +0.44%  "
Set-TextValue "E20" "  +0.08%  "
Set-TextValue "D21" "17.09"
Set-TextValue "E21" "  +0.13%  "
Set-TextValue "D22" "5.746"
Set-TextValue "E22" "  -2.13%  "
Set-TextValue "D23" "27.628.04"
Set-TextValue "E24" "  -0.54%  "
Set-TextValue "D25" "2.090"
Set-TextValue "E25" "  -0.44%  "
Set-TextValue "D26" "165.08"
Set-TextValue "E26" "  +1.53%  "
Set-TextValue "D27" "20.16"
Set-TextValue "E27" "  -1.71%  "
Set-TextValue "D28" "1.959.66"
Set-TextValue "E28" "  -0.05%  "
Set-TextValue "D29" "2.082"
Set-TextValue "E29" "  -3.48%  "
Set-TextValue "D30" "125.39"
Set-TextValue "E30" "  -0.43%  "
Set-TextValue "D31" "1.087"
Set-TextValue "E31" "  -1.39%  "
Set-TextValue "D32" "0.09188"
Set-TextValue "E32" "  +2.54%  "
Set-TextValue "D33" "3.664"
Set-TextValue "E33" "  -0.61%  "
Set-TextValue "D34" "5.502"
Set-TextValue "E34" "  -1.83%  "
Set-TextValue "D35" "0.02287"
Set-TextValue "E35" "  -0.73%  "
Set-TextValue "E36" "  -4.30%  "
Set-TextValue "E37" "  +0.25%  "
Set-TextValue "D38" "0.2090"
Set-TextValue "E38" "  -1.09%  "
Set-TextValue "D39" "0.6285"
Set-TextValue "E39" "  -1.66%  "
Set-TextValue "D40" "4.926"
Set-TextValue "E40" "  -0.50%  "
Set-TextValue "E41" "  -0.62%  "
Set-TextValue "E42" "  -0.33%  "
Set-TextValue "D43" "7.756"
Set-TextValue "E43" "  -2.20%  "
Set-TextValue "D44" "13.18"
Set-TextValue "E44" "  -1.61%  "
Set-TextValue "D45" "3.714"
Set-TextValue "E45" "  -0.01%  "
Set-TextValue "D46" "0.5864"
Set-TextValue "E46" "  -1.10%  "
Set-TextValue "D47" "122.36"
Set-TextValue "E47" "  -0.17%  "
Set-TextValue "D48" "1.939"
Set-TextValue "E48" "  -2.49%  "
Set-TextValue "D49" "0.06885"
Set-TextValue "E49" "  +0.08%  "
Set-TextValue "E50" "  -3.58%  "
Set-TextValue "D51" "71.70"
Set-TextValue "E51" "  -1.93%  "
